$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update first name column (column C) with real user data
$ws.Range("C2").Value = "Aman"
$ws.Range("C3").Value = "Sahil"
$ws.Range("C4").Value = "Rohit"
$ws.Range("C5").Value = "Abhishek"

# Update last name column (column D)
$ws.Range("D2").Value = "Bansal"
$ws.Range("D3").Value = "Awasthi"
$ws.Range("D4").Value = "Dharam"
$ws.Range("D5").Value = "Chaudhary"

# Update pin code column (column E)
$ws.Range("E2").Value = "A12345"
$ws.Range("E3").Value = "A12346"
$ws.Range("E4").Value = "A12347"
$ws.Range("E5").Value = "A12348"

# Update selected cell to match author's final cursor position
$ws.Range("G4").Select()
